# v1.5 modify delete post CRS & review publish audio
# 1/modify delete post CRSs according to reviewer cooments
# 2/review publish audio CRSs, add my comments , ask for some modifications

$wb = $excel.ActiveWorkbook
$review = $wb.Worksheets.Item("REVIEW-SHEET")
$history = $wb.Worksheets.Item("VERSION-HISTORY")

# ---------------------------------------------------------------------
# 1) REVIEW-SHEET: the delete-post review rows (6-10) are now closed
#    (reviewer comments have been addressed)
# ---------------------------------------------------------------------
$review.Range("I6:I10").Value = "closed"

# ---------------------------------------------------------------------
# 2) VERSION-HISTORY: start logging the new v1.5 entry
# ---------------------------------------------------------------------
$history.Range("A7").Value = "v1.5"

# ---------------------------------------------------------------------
# 3) REVIEW-SHEET: three new rows reviewing the "publish audio" CRSs
#    Values are entered column-by-column (matching the order the rows
#    were actually filled in) so new shared-string entries line up.
# ---------------------------------------------------------------------

# Reviewed-entity IDs (column C)
$review.Range("C11").Value = "LH-CRS-PUBLISHAUDIO-001"
$review.Range("C12").Value = "LH-CRS-PUBLISHAUDIO-002"
$review.Range("C13").Value = "LH-CRS-PUBLISHAUDIO-003"

# Review IDs (column B) - row 11 reuses an existing review id
$review.Range("B11").Value = "LH-CRS-Review-008"
$review.Range("B12").Value = "LH-CRS-Review-009"
$review.Range("B13").Value = "LH-CRS-Review-010"

# Version under review (column E)
$review.Range("E11").Value = "v2.2"
$review.Range("E12").Value = "v2.2"
$review.Range("E13").Value = "v2.2"

# Review comments / actions (columns G & F), row by row
$review.Range("G11").Value = "it can be: ""Only registered and logged-in users can access the audio publishing interface that found in a publish drop down in categories page"""
$review.Range("F11").Value = "this CRS contain a lot of details that might be important in SRS not here`nand it can more simple and to the point"

$review.Range("F12").Value = "this CRS contain a lot of details that might be important in SRS not here`nand it can more simple and to the point"
$review.Range("G12").Value = "it can be: ""Only registered and logged-in users can recored a voice note and set a title for it then publish it """

$review.Range("F13").Value = "you mentioned that the audio recorde will have two limits ""one for duration of record and another for its size"" ,but I think we can only suffice with only duration as the 5 minutes audio recorde ussaully doesn't exceeds 20 MB"
$review.Range("G13").Value = "it can be: ""The user can record a voice message for up to 5 minutes. If the time limit is exceeded,a message will appear explaining this."""

# Remaining columns reuse already-known values
$review.Range("A11:A13").Value = "30/4/2025"
$review.Range("D11:D13").Value = "Ahmed Abuzaid"
$review.Range("H11:H13").Value = "Gehad Ashry"
$review.Range("I11:J13").Value = "open"

# Copy the banded-row formatting down from the existing rows (odd/even
# striping: row 9 -> row 11, row 10 -> row 12, row 9 -> row 13) so the
# new rows keep the same look as the rest of the table.
$review.Range("A9:J9").Copy()
$review.Range("A11:J11").PasteSpecial(-4122)
$review.Range("A10:J10").Copy()
$review.Range("A12:J12").PasteSpecial(-4122)
$review.Range("A9:J9").Copy()
$review.Range("A13:J13").PasteSpecial(-4122)

$review.Rows.Item(11).RowHeight = 45
$review.Rows.Item(12).RowHeight = 30
$review.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------------
# 4) VERSION-HISTORY: finish logging the v1.5 entry
# ---------------------------------------------------------------------
$history.Range("B7").Value = "Ahmed Abuzaid"
$history.Range("D7").Value = 45777
$history.Range("C7").Value = "1/modify delete post CRSs according to reviewer cooments`n2/review publish audio CRSs, add my comments , ask for some modifications"

$history.Range("A6:D6").Copy()
$history.Range("A7:D7").PasteSpecial(-4122)
$history.Range("D7").Value = 45777
$history.Rows.Item(7).RowHeight = 93.75

# widen column C on VERSION-HISTORY to fit the longer note
$history.Columns.Item(3).ColumnWidth = 67.67

# ---------------------------------------------------------------------
# 5) View bookkeeping: REVIEW-SHEET becomes the active/selected tab
# ---------------------------------------------------------------------
$history.Range("C9").Select()

$review.Activate()
$review.Range("L13").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
